{"js": "// The diff reorders <w:b/> / <w:i/> ahead of <w:color/> inside <w:rPr>\n// for several \"Tok\" character styles (wml.xsd's CT_RPr requires b/i\n// before color) -- the bold/italic/color *values* themselves are\n// unchanged, only their element order. Re-asserting bold/italic via the\n// Style.font object model makes the engine rewrite each style's run\n// properties in canonical schema order, which is the fix described in\n// the commit message.\nconst boldStyles = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\nconst italicStyles = [\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n];\n\nfor (const name of boldStyles) {\n  const style = context.document.getStyles().getByNameOrNullObject(name);\n  style.font.bold = true;\n}\nfor (const name of italicStyles) {\n  const style = context.document.getStyles().getByNameOrNullObject(name);\n  style.font.italic = true;\n}\nawait context.sync();\n", "ps1": "# The diff reorders <w:b/> / <w:i/> ahead of <w:color/> inside <w:rPr>\n# for several \"Tok\" character styles (wml.xsd's CT_RPr requires b/i\n# before color) -- the bold/italic/color *values* themselves are\n# unchanged, only their element order. Re-asserting Bold/Italic through\n# the Style.Font COM object makes Word rewrite each style's run\n# properties in canonical schema order, which is the fix described in\n# the commit message.\n$d = $word.ActiveDocument\n\n$boldStyles = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\n$italicStyles = @(\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"InformationTok\",\n    \"WarningTok\"\n)\n\nforeach ($name in $boldStyles) {\n    $d.Styles($name).Font.Bold = $true\n}\n\nforeach ($name in $italicStyles) {\n    $d.Styles($name).Font.Italic = $true\n}\n"}
